$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 593.2222
$ws.Range("I9").Value = 593.2222
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 593.2222
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -424.2222
$ws.Range("H34").Value = 3707.8333
$ws.Range("I34").Value = 3707.8333
$ws.Range("K34").Value = 3707.8333
$ws.Range("M34").Value = -3504.8333
$ws.Range("H36").Value = 3707.8333
$ws.Range("I36").Value = 3707.8333
$ws.Range("K36").Value = 3707.8333
$ws.Range("M36").Value = -2992.8333
$ws.Range("H47").Value = 13599.5
$ws.Range("I47").Value = 13599.5
$ws.Range("K47").Value = 13599.5
$ws.Range("M47").Value = -12627.5
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0
$ws.Range("H86").Value = 476143.5
$ws.Range("J86").Value = 476143.5
$ws.Range("L86").Value = 476143.5
$ws.Range("N86").Value = -478389.5
$ws.Range("H87").Value = 116498
$ws.Range("J87").Value = 158996
$ws.Range("L87").Value = 158996
$ws.Range("N87").Value = -161492
$ws.Range("H88").Value = 1731.875
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1731.875
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 1731.875
$ws.Range("N88").Value = -2543.875
$ws.Range("H89").Value = 476143.5
$ws.Range("J89").Value = 476143.5
$ws.Range("L89").Value = 2380717.5
$ws.Range("N89").Value = -2391949.5
$ws.Range("H90").Value = 116498
$ws.Range("J90").Value = 158996
$ws.Range("L90").Value = 476988
$ws.Range("N90").Value = -489468
$ws.Range("H91").Value = 1731.875
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1731.875
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 1731.875
$ws.Range("N91").Value = -4539.875
$ws.Range("H98").Value = 2024.4
$ws.Range("I98").Value = 1975
$ws.Range("K98").Value = 1975
$ws.Range("M98").Value = -477
$ws.Range("H116").Value = 7939.8
$ws.Range("J116").Value = 7939.8
$ws.Range("L116").Value = 7939.8
$ws.Range("N116").Value = -14823.8
$ws.Range("H122").Value = 2024.4
$ws.Range("I122").Value = 1975
$ws.Range("K122").Value = 5925
$ws.Range("M122").Value = -3475
$ws.Range("H131").Value = 15300.286
$ws.Range("I131").Value = 15300.286
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 45900.858
$ws.Range("L131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -40860.858
$ws.Range("H132").Value = 3233.375
$ws.Range("I132").Value = 3233.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9700.125
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -7170.125
$ws.Range("H137").Value = 997.5
$ws.Range("I137").Value = 997.5
$ws.Range("K137").Value = 2992.5
$ws.Range("M137").Value = -442.5
$ws.Range("H138").Value = 1485.6316
$ws.Range("I138").Value = 1485.6316
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4456.8948
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = 683.1052
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4203516.5
$ws.Range("I32").Value = 3892796.2
$ws.Range("K32").Value = 3892796.2
$ws.Range("M32").Value = -3892509.2
$ws.Range("H49").Value = 66000
$ws.Range("J49").Value = 66000
$ws.Range("L49").Value = 66000
$ws.Range("N49").Value = -66520
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 542.2857
$ws.Range("J80").Value = 246.5
$ws.Range("L80").Value = 246.5
$ws.Range("N80").Value = -2242.5
$ws.Range("H83").Value = 542.2857
$ws.Range("J83").Value = 246.5
$ws.Range("L83").Value = 1232.5
$ws.Range("N83").Value = -11216.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1787
$ws.Range("I99").Value = 1338.6
$ws.Range("J99").Value = 2067.25
$ws.Range("K99").Value = 1338.6
$ws.Range("L99").Value = 2067.25
$ws.Range("M99").Value = 159.4000000000001
$ws.Range("N99").Value = -5063.25
$ws.Range("H122").Value = 1695.3334
$ws.Range("I122").Value = 1830.7693
$ws.Range("J122").Value = 1343.2
$ws.Range("K122").Value = 5492.3079
$ws.Range("L122").Value = 4029.6
$ws.Range("M122").Value = -3042.3079
$ws.Range("N122").Value = -8929.6
$ws.Range("H126").Value = 1787
$ws.Range("I126").Value = 1338.6
$ws.Range("J126").Value = 2067.25
$ws.Range("K126").Value = 4015.8
$ws.Range("L126").Value = 6201.75
$ws.Range("M126").Value = -1545.8
$ws.Range("N126").Value = -11141.75
$ws.Range("H134").Value = 2774.2222
$ws.Range("I134").Value = 2746
$ws.Range("K134").Value = 8238
$ws.Range("M134").Value = -5703
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2177.5
$ws.Range("J81").Value = 2177.5
$ws.Range("L81").Value = 6532.5
$ws.Range("N81").Value = -8778.5
$ws.Range("H84").Value = 2177.5
$ws.Range("J84").Value = 2177.5
$ws.Range("L84").Value = 19597.5
$ws.Range("N84").Value = -30829.5
$ws.Range("H122").Value = 150
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 150
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 1350
$ws.Range("N122").Value = -6250
$ws.Range("H136").Value = 2235
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -19200
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3671
$ws.Range("I126").Value = 1739.4
$ws.Range("K126").Value = 5218.200000000001
$ws.Range("M126").Value = -2748.200000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1537.25
$ws.Range("I16").Value = 1537.25
$ws.Range("K16").Value = 1537.25
$ws.Range("M16").Value = -1367.25
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 5000
$ws.Range("N42").Value = -6126
$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("L49").Value = 5000
$ws.Range("N49").Value = -5294
$ws.Range("H82").Value = 1455.875
$ws.Range("J82").Value = 1669.4
$ws.Range("L82").Value = 1669.4
$ws.Range("N82").Value = -2391.4
$ws.Range("H85").Value = 1455.875
$ws.Range("J85").Value = 1669.4
$ws.Range("L85").Value = 1669.4
$ws.Range("N85").Value = -4165.4
$ws.Range("H122").Value = 8001.5
$ws.Range("I122").Value = 8709.5
$ws.Range("J122").Value = 7411.5
$ws.Range("K122").Value = 26128.5
$ws.Range("L122").Value = 22234.5
$ws.Range("M122").Value = -23678.5
$ws.Range("N122").Value = -27134.5
$ws.Range("H132").Value = 3256
$ws.Range("I132").Value = 3075.1667
$ws.Range("K132").Value = 9225.500100000001
$ws.Range("M132").Value = -6695.500100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2371.818
$ws.Range("I126").Value = 2254.6667
$ws.Range("K126").Value = 6764.000100000001
$ws.Range("M126").Value = -4294.000100000001

Write-Output "Applied 197 cell changes across 8 sheets"
